# Update Name of Algo
# Apply updated numeric results to the worksheet (RandomForest imputation output)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = -11.73719999999999
$ws.Range("B3").Value = 6.080499999999987
$ws.Range("C5").Value = -14.0278
$ws.Range("E5").Value = 12.45879999999999
$ws.Range("E9").Value = 13.77560000000001
$ws.Range("E11").Value = 13.6764
$ws.Range("B14").Value = 9.322300000000006
$ws.Range("B21").Value = 5.908099999999994
$ws.Range("E21").Value = 13.20459999999999
$ws.Range("B23").Value = 5.632200000000001
$ws.Range("B25").Value = 5.842799999999994
